# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" stat for each row. Update the computed values
# for rows 2-15 to reflect the regenerated K values (Strike# replaced by K).
$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 2
    6  = 0
    7  = 2
    8  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 2
    15 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
